# data: add 2023-05-29 notices
#
# Inserts the 宏昌科技 (301008) notice as a new row 3 on both sheets,
# pushing the existing 祥源新材 / 佳禾智能 rows down by one. Also tidies up
# sheet 1's formatting to match sheet 2 (bold/boxed header style, no more
# hyperlink styling/relationships on column F, default page margins).

$wb  = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("公告明细")
$ws2 = $wb.Worksheets.Item("公告汇总")

# ---------------------------------------------------------------------
# Sheet 1: 公告明细 (notice detail)
# ---------------------------------------------------------------------

# Insert a new row 3 (祥源新材 -> row4, 佳禾智能 -> row5) and fill it in.
$ws1.Rows.Item(3).Insert()

$ws1.Range("A3").Value = "宏昌科技"
$ws1.Range("B3").Value = "'301008"
$ws1.Range("C3").Value = "宏昌科技:关于使用部分闲置募集资金及自有资金进行现金管理的进展公告"
$ws1.Range("D3").Value = "2023-05-30 00:00:00"
$ws1.Range("E3").Value = "2023-05-29 16:28:45:000"
$ws1.Range("F3").Value = "https://pdf.dfcfw.com/pdf/H2_AN202305291587334938_1.pdf?1685377738000.pdf"
# Stock code is textual (keep leading format), not a number.
$ws1.Range("B3").Style = "Normal"

# Column F is no longer a set of live hyperlinks - drop them and their style.
$ws1.Cells.Hyperlinks.Delete()
$ws1.Range("F2:F5").Style = "Normal"

# Match the bold / bordered header formatting already used on sheet 2.
$ws2.Range("A1").Copy()
$ws1.Range("A1:F1").PasteSpecial(-4122)

# Reset the page margins to Excel's regular defaults.
$ws1.PageSetup.LeftMargin = 54
$ws1.PageSetup.RightMargin = 54
$ws1.PageSetup.TopMargin = 72
$ws1.PageSetup.BottomMargin = 72
$ws1.PageSetup.HeaderMargin = 36
$ws1.PageSetup.FooterMargin = 36

# ---------------------------------------------------------------------
# Sheet 2: 公告汇总 (notice summary)
# ---------------------------------------------------------------------

# Insert a new row 3 (祥源新材 -> row4, 佳禾智能 -> row5) and fill it in.
$ws2.Rows.Item(3).Insert()

$ws2.Range("A3").Value = "宏昌科技"
$ws2.Range("B3").Value = "'301008"
$ws2.Range("C3").Value = 1
$ws2.Range("B3").Style = "Normal"
